$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Type BUY/SELL (C5) and Amount (D5)
$ws.Range("C5").Value = "BUY"
$ws.Range("D5").Value = 33.5321

# Row 18: Amount (D18)
$ws.Range("D18").Value = 8.2939000000000007

# Row 28: Amount (D28)
$ws.Range("D28").Value = 3.7

# Update selected cell to mirror sheetView selection change
$ws.Range("H16").Select()
